# The commit replaces the text in cell E8 ("Good Morning" -> "GIT UPDATE")
# and leaves that cell selected, matching the author's last interactive edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
